$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Foglio1 (sheet1)
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Foglio1")

# Latenza post-memoria count corrected from 5 to 6
$ws1.Range("A55").Value = 6

# Clear the old "Esempio 16x16 / 64x64 / 32x32" blocks (rows 56-76) so we
# can rebuild them in the new layout (the "scrittura nel rf di uscita" row
# is removed, the 32x32 block now precedes the 64x64 block, and the post
# memory latency counts are corrected from 5 to 6).
$ws1.Range("A56:R76").ClearContents()

# --- Esempio 16x16 block (rows 50-56) ---
$ws1.Range("A56").Formula = "=SUM(A51:A55)"

# --- Esempio 32x32 block (rows 58-64) ---
$ws1.Range("A58").Value = "Esempio 32x32"
$ws1.Range("A59").Value = 8
$ws1.Range("B59").Value = "Register_in"
$ws1.Range("A60").Value = 64
$ws1.Range("B60").Value = "elaborazione dati (4 blocchi 16x16)"
$ws1.Range("A61").Value = 0
$ws1.Range("B61").Value = "blocchi persi a causa del costruttore"
$ws1.Range("A62").Value = 64
$ws1.Range("B62").Value = "elaborazione dati (4 blocchi 16x16)"
$ws1.Range("A63").Value = 6
$ws1.Range("B63").Value = "Latenza post-memoria (scelta del best candidate)"
$ws1.Range("A64").Formula = "=SUM(A59:A63)"

# --- Esempio 64x64 block (rows 66-72) ---
$ws1.Range("A66").Value = "Esempio 64x64"
$ws1.Range("A67").Value = 8
$ws1.Range("B67").Value = "Register_in"
$ws1.Range("A68").Value = 256
$ws1.Range("B68").Value = "elaborazione dati (16 blocchi 16x16)"
$ws1.Range("A69").Value = 0
$ws1.Range("B69").Value = "blocchi persi a causa del costruttore"
$ws1.Range("A70").Value = 256
$ws1.Range("B70").Value = "elaborazione dati (16 blocchi 16x16)"
$ws1.Range("A71").Value = 6
$ws1.Range("B71").Value = "Latenza post-memoria (scelta del best candidate)"
$ws1.Range("A72").Formula = "=SUM(A67:A71)"

# CLK-cycle totals now reference the (corrected) computed sums instead of
# hard-coded constants.
$ws1.Range("N11").Formula = "=A56"
$ws1.Range("N13").Formula = "=A64"
$ws1.Range("N15").Formula = "=A72"

$ws1.Range("L20").Select() | Out-Null

# -----------------------------------------------------------------
# Foglio2 (sheet2)
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Foglio2")

$ws2.Range("J11").Formula = "=A56"
$ws2.Range("J13").Formula = "=A64"
$ws2.Range("J15").Formula = "=A72"

# --- Esempio 16x16 block (rows 50-56) ---
$ws2.Range("A50").Value = "Esempio 16x16"
$ws2.Range("A51").Value = 8
$ws2.Range("B51").Value = "Register_in"
$ws2.Range("A52").Value = 16
$ws2.Range("B52").Value = "Parallelism"
$ws2.Range("A53").Value = 4
$ws2.Range("B53").Value = "11 to 27 b"
$ws2.Range("A54").Value = 16
$ws2.Range("B54").Value = "Parallelism"
$ws2.Range("A55").Value = 6
$ws2.Range("B55").Value = "Latenza post-memoria (scelta del best candidate)"
$ws2.Range("A56").Formula = "=SUM(A51:A55)"

# --- Esempio 32x32 block (rows 58-64) ---
$ws2.Range("A58").Value = "Esempio 32x32"
$ws2.Range("A59").Value = 8
$ws2.Range("B59").Value = "Register_in"
$ws2.Range("A60").Value = 64
$ws2.Range("B60").Value = "elaborazione dati (4 blocchi 16x16)"
$ws2.Range("A61").Value = 0
$ws2.Range("B61").Value = "blocchi persi a causa del costruttore"
$ws2.Range("A62").Value = 64
$ws2.Range("B62").Value = "elaborazione dati (4 blocchi 16x16)"
$ws2.Range("A63").Value = 6
$ws2.Range("B63").Value = "Latenza post-memoria (scelta del best candidate)"
$ws2.Range("A64").Formula = "=SUM(A59:A63)"

# --- Esempio 64x64 block (rows 66-72) ---
$ws2.Range("A66").Value = "Esempio 64x64"
$ws2.Range("A67").Value = 8
$ws2.Range("B67").Value = "Register_in"
$ws2.Range("A68").Value = 256
$ws2.Range("B68").Value = "elaborazione dati (16 blocchi 16x16)"
$ws2.Range("A69").Value = 0
$ws2.Range("B69").Value = "blocchi persi a causa del costruttore"
$ws2.Range("A70").Value = 256
$ws2.Range("B70").Value = "elaborazione dati (16 blocchi 16x16)"
$ws2.Range("A71").Value = 6
$ws2.Range("B71").Value = "Latenza post-memoria (scelta del best candidate)"
$ws2.Range("A72").Formula = "=SUM(A67:A71)"

$ws2.Range("J16").Select() | Out-Null

# -----------------------------------------------------------------
# Foglio3 (sheet3)
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Foglio3")

$ws3.Range("L11").Value = 50
$ws3.Range("L13").Value = 142
$ws3.Range("L15").Value = 526

$ws3.Range("K26").Select() | Out-Null

$wb.Application.Calculate()
